$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 10).Value = 5
$ws.Cells.Item(9, 10).Value = 10
$ws.Cells.Item(10, 10).Value = 15
$ws.Cells.Item(15, 10).Value = 2
$ws.Cells.Item(16, 10).Value = 2

$ws.Range("K15").Select()
